# Data Type.docx edit script
# 1. Remove the stray "_GoBack" bookmark paragraph that sits right after
#    "Length of stay is an integer greater than 0."
# 2. Prefix the first Date/Time cells (check-in/out table) with "c"
#    as a separate run.
# 3. Prefix the second Date/Time cells (payment table) with "p"
#    as a separate run.
# 4. Split the second "Bill_event_payment" cell into "check_out" +
#    a relocated "_GoBack" bookmark + "_payment_id".

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------
# Step 1: remove the old "_GoBack" bookmark paragraph.
# ---------------------------------------------------------------------
$bmRng = $d.Content
$bmRng.Find.Execute("Length of stay is an integer greater than 0.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmParaRng = $bmRng.Next(4, 1)
$bmParaRng.InsertXML("<w:p $wNs/>")

# ---------------------------------------------------------------------
# Helper: prefix the lone run of the Nth paragraph whose full text
# equals $label with a new run containing $prefix, preserving the
# paragraph's own attributes / pPr.
# ---------------------------------------------------------------------
function Add-PrefixRun($label, $prefix, $occurrence) {
    $script:xml = $d.Content.WordOpenXML
    $pattern = '<w:p ([^>]*)>(?:(<w:pPr>.*?</w:pPr>))?<w:r><w:t>' + $label + '</w:t></w:r></w:p>'
    $re = [regex]$pattern
    $ms = $re.Matches($script:xml)
    $m = $ms.Item($occurrence - 1)
    $attrs = $m.Groups[1].Value
    $ppr = $m.Groups[2].Value

    $rng = $d.Content
    for ($i = 1; $i -le $occurrence; $i++) {
        if ($i -gt 1) { $rng.Collapse(0) }
        $rng.Find.Execute($label, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    }
    $s = $rng.Start
    $e = $rng.End
    $target = $d.Range($s, $e)
    $newXml = "<w:p $wNs $w14Ns $attrs>$ppr<w:r><w:t>$prefix</w:t></w:r><w:r><w:t>$label</w:t></w:r></w:p>"
    $target.InsertXML($newXml)
}

Add-PrefixRun "Date" "c" 1
Add-PrefixRun "Time" "c" 1
Add-PrefixRun "Date" "p" 2
Add-PrefixRun "Time" "p" 2

# ---------------------------------------------------------------------
# Step 6: split the second "Bill_event_payment" into
#         "check_out" + bookmark(_GoBack) + "_payment_id"
# ---------------------------------------------------------------------
$xml = $d.Content.WordOpenXML
$bepRe = [regex]'<w:p ([^>]*)><w:pPr>(<w:jc w:val="center"/>)</w:pPr><w:r><w:t>Bill_event_payment</w:t></w:r></w:p>'
$bepMs = $bepRe.Matches($xml)
$bepM = $bepMs.Item($bepMs.Count - 1)
$bepAttrs = $bepM.Groups[1].Value
$bepPpr = $bepM.Groups[2].Value

$rng = $d.Content
$rng.Find.Execute("Bill_event_payment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.Find.Execute("Bill_event_payment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $rng.Start
$e = $rng.End
$target = $d.Range($s, $e)
$newXml = "<w:p $wNs $w14Ns $bepAttrs><w:pPr>$bepPpr</w:pPr><w:r><w:t>check_out</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t>_payment_id</w:t></w:r></w:p>"
$target.InsertXML($newXml)

Write-Output "All edits applied."
